# Add a "Note: ..." paragraph (styled BodyText) right after the paragraph
# that ends the ELISA procedure description ("... subsequently repeated).")
# and before the "References" heading.

$d = $word.ActiveDocument

# Locate the last paragraph of the assay-description section using a short,
# unique, accent-free anchor so the Find is not sensitive to the exact
# (non-breaking) spaces/µ characters used in the surrounding text.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "(that are subsequently repeated).",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

# Collapse to the end of the match (end of that paragraph's text) and split
# off a brand-new, empty paragraph right after it.
$anchor.Collapse(0)
$anchor.InsertParagraphAfter() | Out-Null

# The freshly created paragraph is the one immediately following the anchor
# paragraph; give it the "BodyText" style, matching the target markup.
$newParaIndex = $anchor.Paragraphs.First.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.Style = "BodyText"
$insertStart = $newPara.Range.Start

# Insert all the plain text first (each call appends right after the
# previous one), then go back and apply character formatting to the
# specific sub-ranges. Setting Font properties directly on a collapsed
# (zero-length) range before inserting text can bleed formatting into
# unrelated, later content, so we avoid that pattern entirely.
$segments = @(
    "Note",
    ": absorbance values need to be written within the Excel template provided by the",
    " ",
    "insane",
    " ",
    "package."
)

$cursor = $insertStart
$bounds = @{}
foreach ($seg in $segments) {
    $segStart = $cursor
    $run = $d.Range($segStart, $segStart)
    $run.InsertAfter($seg)
    $cursor = $segStart + $seg.Length
    $bounds[$seg] = @($segStart, $cursor)
}

# "Note" -> italic
$noteBounds = $bounds["Note"]
$noteRange = $d.Range($noteBounds[0], $noteBounds[1])
$noteRange.Font.Italic = $true

# "insane" -> bold
$insaneBounds = $bounds["insane"]
$insaneRange = $d.Range($insaneBounds[0], $insaneBounds[1])
$insaneRange.Font.Bold = $true

Write-Output "Inserted note paragraph after anchor; new paragraph index=$newParaIndex"
